# Update TPM-derived NATMI ligand-receptor metrics (Sema7a-Itgb1) on Sheet1.
# Columns: G/H = ligand avg/total expr, I/J = ligand specificity,
#          M/N = receptor avg/total expr, O/P = receptor specificity,
#          Q/R = edge avg/total weight, S/T = edge specificity.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> Sema7a -> Itgb1 -> ECs)
$ws.Range("G2").Value = 20.09599933333334
$ws.Range("H2").Value = 60.287998
$ws.Range("I2").Value = 0.9230842418515092
$ws.Range("J2").Value = 0.9230842418515091
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 1226.692080180263
$ws.Range("R2").Value = 11040.22872162236
$ws.Range("S2").Value = 0.1886427381891983
$ws.Range("T2").Value = 0.1886427381891983

# Row 3 (ECs -> Sema7a -> Itgb1 -> FAPs)
$ws.Range("G3").Value = 20.09599933333334
$ws.Range("H3").Value = 60.287998
$ws.Range("I3").Value = 0.9230842418515092
$ws.Range("J3").Value = 0.9230842418515091
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("Q3").Value = 2136.495437859689
$ws.Range("R3").Value = 19228.4589407372
$ws.Range("S3").Value = 0.3285538042010966
$ws.Range("T3").Value = 0.3285538042010966

# Row 4 (ECs -> Sema7a -> Itgb1 -> MuSCs)
$ws.Range("G4").Value = 20.09599933333334
$ws.Range("H4").Value = 60.287998
$ws.Range("I4").Value = 0.9230842418515092
$ws.Range("J4").Value = 0.9230842418515091
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 2639.376586403727
$ws.Range("R4").Value = 23754.38927763355
$ws.Range("S4").Value = 0.4058876994612142
$ws.Range("T4").Value = 0.4058876994612142

# Row 5 (FAPs -> Sema7a -> Itgb1 -> ECs)
$ws.Range("I5").Value = 0.05319611498621682
$ws.Range("J5").Value = 0.05319611498621682
$ws.Range("M5").Value = 61.04160633333334
$ws.Range("N5").Value = 183.124819
$ws.Range("O5").Value = 0.2043613460574534
$ws.Range("P5").Value = 0.2043613460574534
$ws.Range("Q5").Value = 70.69263019706924
$ws.Range("R5").Value = 636.2336717736231
$ws.Range("S5").Value = 0.01087122966361034
$ws.Range("T5").Value = 0.01087122966361034

# Row 6 (FAPs -> Sema7a -> Itgb1 -> FAPs)
$ws.Range("I6").Value = 0.05319611498621682
$ws.Range("J6").Value = 0.05319611498621682
$ws.Range("O6").Value = 0.3559304658284363
$ws.Range("P6").Value = 0.3559304658284363
$ws.Range("S6").Value = 0.01893411798730721
$ws.Range("T6").Value = 0.01893411798730722

# Row 7 (FAPs -> Sema7a -> Itgb1 -> MuSCs)
$ws.Range("I7").Value = 0.05319611498621682
$ws.Range("J7").Value = 0.05319611498621682
$ws.Range("M7").Value = 131.3384093333333
$ws.Range("N7").Value = 394.015228
$ws.Range("O7").Value = 0.4397081881141102
$ws.Range("P7").Value = 0.4397081881141103
$ws.Range("Q7").Value = 152.1037560999195
$ws.Range("R7").Value = 1368.933804899276
$ws.Range("S7").Value = 0.02339076733529926
$ws.Range("T7").Value = 0.02339076733529927

# Row 8 (MuSCs -> Sema7a -> Itgb1 -> ECs)
$ws.Range("G8").Value = 0.5163883333333333
$ws.Range("H8").Value = 1.549165
$ws.Range("I8").Value = 0.02371964316227407
$ws.Range("J8").Value = 0.02371964316227407
$ws.Range("M8").Value = 61.04160633333334
$ws.Range("N8").Value = 183.124819
$ws.Range("O8").Value = 0.2043613460574534
$ws.Range("P8").Value = 0.2043613460574534
$ws.Range("Q8").Value = 31.52117335845945
$ws.Range("R8").Value = 283.690560226135
$ws.Range("S8").Value = 0.004847378204644801
$ws.Range("T8").Value = 0.004847378204644801

# Row 9 (MuSCs -> Sema7a -> Itgb1 -> FAPs)
$ws.Range("G9").Value = 0.5163883333333333
$ws.Range("H9").Value = 1.549165
$ws.Range("I9").Value = 0.02371964316227407
$ws.Range("J9").Value = 0.02371964316227407
$ws.Range("O9").Value = 0.3559304658284363
$ws.Range("P9").Value = 0.3559304658284363
$ws.Range("Q9").Value = 54.89954990696333
$ws.Range("R9").Value = 494.09594916267
$ws.Range("S9").Value = 0.008442543640032494
$ws.Range("T9").Value = 0.008442543640032496

# Row 10 (MuSCs -> Sema7a -> Itgb1 -> MuSCs)
$ws.Range("G10").Value = 0.5163883333333333
$ws.Range("H10").Value = 1.549165
$ws.Range("I10").Value = 0.02371964316227407
$ws.Range("J10").Value = 0.02371964316227407
$ws.Range("M10").Value = 131.3384093333333
$ws.Range("N10").Value = 394.015228
$ws.Range("O10").Value = 0.4397081881141102
$ws.Range("P10").Value = 0.4397081881141103
$ws.Range("Q10").Value = 67.8216222982911
$ws.Range("R10").Value = 610.3946006846199
$ws.Range("S10").Value = 0.01042972131759678
$ws.Range("T10").Value = 0.01042972131759678
